$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 'dnasr281@gmail.com, System'
    3 = 'dnasr281@gmail.com, System'
    4 = 'dnasr281@gmail.com, System'
    5 = 'dnasr281@gmail.com, System'
    6 = 'dnasr281@gmail.com, System'
    7 = 'dnasr281@gmail.com, System'
    8 = 'System, dnasr281@gmail.com'
    16 = 'dnasr281@gmail.com, System'
    17 = 'dnasr281@gmail.com, System'
    22 = 'dnasr281@gmail.com, System'
    23 = 'dnasr281@gmail.com, System'
    24 = 'System, dnasr281@gmail.com'
    26 = 'System, dnasr281@gmail.com'
    29 = 'System, dnasr281@gmail.com'
    32 = 'System, dnasr281@gmail.com'
    37 = 'dnasr281@gmail.com, System'
    38 = 'dnasr281@gmail.com, System'
    43 = 'dnasr281@gmail.com, System'
    44 = 'dnasr281@gmail.com, System'
    45 = 'System, dnasr281@gmail.com'
    47 = 'System, dnasr281@gmail.com'
    50 = 'System, dnasr281@gmail.com'
    53 = 'System, dnasr281@gmail.com'
    58 = 'dnasr281@gmail.com, System'
    59 = 'dnasr281@gmail.com, System'
    64 = 'dnasr281@gmail.com, System'
    65 = 'dnasr281@gmail.com, System'
    66 = 'System, dnasr281@gmail.com'
    68 = 'System, dnasr281@gmail.com'
    71 = 'System, dnasr281@gmail.com'
    74 = 'System, dnasr281@gmail.com'
    79 = 'dnasr281@gmail.com, System'
    80 = 'dnasr281@gmail.com, System'
    85 = 'dnasr281@gmail.com, System'
    86 = 'dnasr281@gmail.com, System'
    87 = 'dnasr281@gmail.com, System'
    88 = 'dnasr281@gmail.com, System'
    89 = 'dnasr281@gmail.com, System'
    90 = 'dnasr281@gmail.com, System'
    91 = 'System, dnasr281@gmail.com'
    99 = 'dnasr281@gmail.com, System'
    100 = 'dnasr281@gmail.com, System'
    105 = 'dnasr281@gmail.com, System'
    106 = 'dnasr281@gmail.com, System'
    107 = 'dnasr281@gmail.com, System'
    108 = 'dnasr281@gmail.com, System'
    109 = 'dnasr281@gmail.com, System'
    110 = 'dnasr281@gmail.com, System'
    111 = 'System, dnasr281@gmail.com'
    119 = 'dnasr281@gmail.com, System'
    120 = 'dnasr281@gmail.com, System'
    125 = 'dnasr281@gmail.com, System'
    126 = 'dnasr281@gmail.com, System'
    127 = 'dnasr281@gmail.com, System'
    128 = 'dnasr281@gmail.com, System'
    129 = 'dnasr281@gmail.com, System'
    130 = 'dnasr281@gmail.com, System'
    131 = 'System, dnasr281@gmail.com'
    139 = 'dnasr281@gmail.com, System'
    140 = 'dnasr281@gmail.com, System'
    145 = 'dnasr281@gmail.com, System'
    146 = 'dnasr281@gmail.com, System'
    147 = 'dnasr281@gmail.com, System'
    148 = 'dnasr281@gmail.com, System'
    149 = 'dnasr281@gmail.com, System'
    150 = 'dnasr281@gmail.com, System'
    151 = 'System, dnasr281@gmail.com'
    159 = 'dnasr281@gmail.com, System'
    160 = 'dnasr281@gmail.com, System'
    165 = 'dnasr281@gmail.com, System'
    166 = 'dnasr281@gmail.com, System'
    167 = 'dnasr281@gmail.com, System'
    168 = 'dnasr281@gmail.com, System'
    169 = 'dnasr281@gmail.com, System'
    170 = 'dnasr281@gmail.com, System'
    171 = 'System, dnasr281@gmail.com'
    179 = 'dnasr281@gmail.com, System'
    180 = 'dnasr281@gmail.com, System'
    185 = 'dnasr281@gmail.com, System'
    186 = 'dnasr281@gmail.com, System'
    187 = 'System, dnasr281@gmail.com'
    189 = 'System, dnasr281@gmail.com'
    192 = 'System, dnasr281@gmail.com'
    195 = 'System, dnasr281@gmail.com'
    200 = 'dnasr281@gmail.com, System'
    201 = 'dnasr281@gmail.com, System'
    206 = 'dnasr281@gmail.com, System'
    207 = 'dnasr281@gmail.com, System'
    208 = 'System, dnasr281@gmail.com'
    210 = 'System, dnasr281@gmail.com'
    213 = 'System, dnasr281@gmail.com'
    216 = 'System, dnasr281@gmail.com'
    221 = 'dnasr281@gmail.com, System'
    222 = 'dnasr281@gmail.com, System'
    227 = 'dnasr281@gmail.com, System'
    228 = 'dnasr281@gmail.com, System'
    229 = 'System, dnasr281@gmail.com'
    231 = 'System, dnasr281@gmail.com'
    234 = 'System, dnasr281@gmail.com'
    237 = 'System, dnasr281@gmail.com'
    242 = 'dnasr281@gmail.com, System'
    243 = 'dnasr281@gmail.com, System'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $values[$row]
}
